$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (Changed) date column (C) for rows 2-11
# from serial date 45202 (2023-10-03) to 45203 (2023-10-04)
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value = 45203
    }
}
